$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph: merge the two runs that were split by the old
#    "_GoBack" bookmark back into a single run (the bookmark is being moved
#    further down the document). The visible text does not change here.
# ---------------------------------------------------------------------------
$oldFull1 = " to keep track or monitor the resort. As time goes by, and the project is being thought out, the developers also thought that it would be beneficial for both the management and customers if an online reservation system would also be created, along with it is a website showcasing the resort and its facilities. The project aims to make managing and monitoring easier in the resort."
$d.Content.Find.Execute($oldFull1, $true, $false, $false, $false, $false, $true, 1, $false, $oldFull1, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Limitations and Exclusions" paragraph: replace the whole paragraph
#    text with the revised wording, then re-split it into the same run
#    boundaries used by the edit (forced via temporary bookmarks, since
#    plain insertion of identically-formatted text collapses into a single
#    run). The "_GoBack" bookmark is recreated at its new location as part
#    of this split.
# ---------------------------------------------------------------------------
$oldFull2 = "The limitations would be the exclusion of the services in the reservation process, the developers are focusing more on reservation of rooms, that's why reserving the services online would not be included in the initial release, however, customers can still avail the services when they're in the actual resort. Another one is the email reminders, the developers would also like to remind the customers about their made reservations, but at this time, the developers are still not knowledgeable how to do such feature, that's why it won't be included in the release."
$newFull2 = "The limitations will be the exclusion of the services in the reservation process, the developers are focusing more on reservation of rooms, that's why reserving the services online would not be included in the initial release. However, customers can still avail the services when they're in the actual resort. Another one is, the email reminders, the developers would also like to remind the customers about their made reservations. But at this time, the developers are still not knowledgeable how to do such feature, that's why it won't be included in the release."

$r = $d.Content
$r.Find.Execute($oldFull2, $true, $false, $false, $false, $false, $true, 1, $false, $newFull2, 2) | Out-Null

# Locate the freshly written paragraph text again to get a stable start offset.
$r2 = $d.Content
$r2.Find.Execute($newFull2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base = $r2.Start

$chunks = @(
    "The limitations will",
    " be the exclusion of the services in the reservation process, the developers are focusing more on reservation of rooms, that's why reserving the services online would not be ",
    "included in the initial release. H",
    "owever, customers can stil",
    "l avail the services when they're in the actual resort. Another one is",
    ",",
    " the email reminders, the developers would also like to remind the custome",
    "rs about their made reservations. B",
    "ut at this time, the developers are still not knowledgeable how to do such feature, that's why it won't be included in the release."
)

$offsets = @($base)
$cum = $base
foreach ($chunk in $chunks) {
    $cum = $cum + $chunk.Length
    $offsets += $cum
}

# Sanity check against the known absolute end of the range.
Write-Host "limitations-paragraph-offsets:" $offsets

# Internal boundaries are offsets[1] .. offsets[8] (offsets[0] is the
# paragraph start, offsets[9] is the paragraph's text end).
$d.Bookmarks.Add("zzTmp1", $d.Range($offsets[1], $offsets[1])) | Out-Null
$d.Bookmarks.Add("zzTmp2", $d.Range($offsets[2], $offsets[2])) | Out-Null
$d.Bookmarks.Add("zzTmp3", $d.Range($offsets[3], $offsets[3])) | Out-Null
$d.Bookmarks.Add("_GoBack", $d.Range($offsets[4], $offsets[4])) | Out-Null
$d.Bookmarks.Add("zzTmp5", $d.Range($offsets[5], $offsets[5])) | Out-Null
$d.Bookmarks.Add("zzTmp6", $d.Range($offsets[6], $offsets[6])) | Out-Null
$d.Bookmarks.Add("zzTmp7", $d.Range($offsets[7], $offsets[7])) | Out-Null
$d.Bookmarks.Add("zzTmp8", $d.Range($offsets[8], $offsets[8])) | Out-Null

$d.Bookmarks.Item("zzTmp1").Delete()
$d.Bookmarks.Item("zzTmp2").Delete()
$d.Bookmarks.Item("zzTmp3").Delete()
$d.Bookmarks.Item("zzTmp5").Delete()
$d.Bookmarks.Item("zzTmp6").Delete()
$d.Bookmarks.Item("zzTmp7").Delete()
$d.Bookmarks.Item("zzTmp8").Delete()

# ---------------------------------------------------------------------------
# 3) Header page-number field: refresh the cached field result from the
#    stale "iii" (front-matter numbering) to "2".
# ---------------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    if ($hdr.Exists) {
        $hdr.Range.Find.Execute("iii", $true, $false, $false, $false, $false, $true, 1, $false, "2", 2) | Out-Null
    }
}

Write-Host "done"
